$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows to match repulled/recalculated data
$ws.Range("F12").Value = -1
$ws.Range("F16").Value = 1
$ws.Range("F22").Value = -4
$ws.Range("F25").Value = -3
$ws.Range("F27").Value = 1
$ws.Range("F33").Value = -5
